$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.436.33'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.856.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.42%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6338'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.98%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07593'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.12%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2926'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.30%  '

$ws.Range("E10").Value = '  -0.91%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07758'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.856.56'
$ws.Range("D12").Style = "Normal"

$ws.Range("E13").Value = '  +0.11%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6863'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.82%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001047'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.40'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.114.38'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.13%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.155'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.458.88'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.01%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '230.41'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.96%  '

$ws.Range("E21").Value = '  -0.02%  '

$ws.Range("E22").Value = '  +0.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.540'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.63%  '

$ws.Range("E24").Value = '  +0.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.15'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.68%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1400'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.79%  '

$ws.Range("E27").Value = '  +0.97%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.76'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.416'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.480'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.10%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05704'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.160'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.88%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.065'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.98%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.833'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.51%  '

$ws.Range("E35").Value = '  -0.32%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6993'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.39%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.250.91'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.83%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01829'
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.776'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.524'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9118'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.92%  '

$ws.Range("E43").Value = '  +0.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.017.16'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.89%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.52'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.20%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.07%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.168'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1168'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.073'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3974'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.678'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.22%  '

